$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-121 down to 45-122.
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = 'Vega Modelo de Temuco'
$ws.Range("C44").Value = 'La Araucanía'
$ws.Range("D44").Value = 44526
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100114007
$ws.Range("G44").Value = 'Jengibre'
$ws.Range("H44").Value = 'Sin especificar'
$ws.Range("I44").Value = 'Primera'
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 20000
$ws.Range("N44").Value = '$/caja 13 kilos'
$ws.Range("O44").Value = 'Perú'
$ws.Range("P44").Value = 1538
$ws.Range("Q44").Value = 13
$ws.Range("R44").Value = 'Hortaliza'
